$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.968.82'
$ws.Range("E2").Value = '  +0.14%  '
$ws.Range("D3").Value = '1.561.56'
$ws.Range("E3").Value = '  +0.45%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = "'" + '207.38'
$ws.Range("E5").Value = '  +0.40%  '
$ws.Range("E6").Value = '  +0.38%  '
$ws.Range("E7").Value = '  +0.11%  '
$ws.Range("E8").Value = '  +0.86%  '
$ws.Range("E9").Value = '  +0.55%  '
$ws.Range("E10").Value = '  +2.69%  '
$ws.Range("E11").Value = '  -0.28%  '
$ws.Range("D12").Value = '1.784.54'
$ws.Range("E12").Value = '  +0.50%  '
$ws.Range("D13").Value = '1.563.04'
$ws.Range("E13").Value = '  +0.56%  '
$ws.Range("E14").Value = '  +0.45%  '
$ws.Range("D15").Value = "'" + '0.520'
$ws.Range("E15").Value = '  +0.72%  '
$ws.Range("D16").Value = "'" + '61.90'
$ws.Range("E16").Value = '  +0.38%  '
$ws.Range("D17").Value = '26.959.23'
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("E18").Value = '  +1.83%  '
$ws.Range("D19").Value = "'" + '215.62'
$ws.Range("E19").Value = '  -1.12%  '
$ws.Range("D20").Value = "'" + '7.35'
$ws.Range("E20").Value = '  +0.95%  '
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("E22").Value = '  +1.53%  '
$ws.Range("E23").Value = '  +0.11%  '
$ws.Range("E24").Value = '  -0.95%  '
$ws.Range("D25").Value = "'" + '153.32'
$ws.Range("E25").Value = '  -0.48%  '
$ws.Range("E26").Value = '  +0.41%  '
$ws.Range("D27").Value = "'" + '15.08'
$ws.Range("E27").Value = '  +1.17%  '
$ws.Range("E28").Value = '  +1.40%  '
$ws.Range("E29").Value = '  +0.01%  '
$ws.Range("E30").Value = '  +0.58%  '
$ws.Range("D31").Value = "'" + '1.11'
$ws.Range("E31").Value = '  +1.57%  '
$ws.Range("E32").Value = '  +0.33%  '
$ws.Range("E33").Value = '  +1.46%  '
$ws.Range("D34").Value = '1.421.28'
$ws.Range("E34").Value = '  -0.84%  '
$ws.Range("E35").Value = '  +2.33%  '
$ws.Range("D36").Value = "'" + '1.06'
$ws.Range("E36").Value = '  +8.60%  '
$ws.Range("D37").Value = "'" + '2.34'
$ws.Range("E37").Value = '  +2.27%  '
$ws.Range("E38").Value = '  +0.00%  '
$ws.Range("D39").Value = "'" + '0.532'
$ws.Range("E39").Value = '  +2.29%  '
$ws.Range("D40").Value = "'" + '5.82'
$ws.Range("E40").Value = '  +2.21%  '
$ws.Range("E41").Value = '  -0.27%  '
$ws.Range("E42").Value = '  +0.06%  '
$ws.Range("E43").Value = '  +2.54%  '
$ws.Range("E44").Value = '  +2.11%  '
$ws.Range("D45").Value = "'" + '64.55'
$ws.Range("E45").Value = '  +1.05%  '
$ws.Range("E46").Value = '  -0.97%  '
$ws.Range("D47").Value = '1.698.10'
$ws.Range("E48").Value = '  +0.23%  '
$ws.Range("E49").Value = '  -0.67%  '
$ws.Range("D50").Value = '0.0₆0100'
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("E51").Value = '  +0.50%  '
